# Updated Indonesia files compatible with v3.3.1
# Target workbook: InputData/bldgs/EoBSDwEC/Elast of Bldg Svc Demand wrt E Cost.xlsx
#
# Summary of the data changes made on the "EoBSDwEC" sheet:
#   - A1 header changes from "Fuel" to "Elasticity by Fuel (dimensionless)",
#     gets word-wrap turned on and the row is made taller to fit two lines.
#   - Rows 2-6 (electricity, coal, natural gas, petroleum diesel, heat) are
#     given an explicit (slightly taller) row height.
#   - Four new fuel rows are appended after "biomass": kerosene, heavy or
#     residual fuel oil, LPG propane or butane, and hydrogen - each carrying
#     the same elasticity values as the other non-electricity fuels
#     (-0.15 / -0.15 / -0.25).
#   - Column widths are tweaked slightly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EoBSDwEC")

# --- New fuel rows appended after "biomass" (row 7) ---------------------
# (written before the A1 header text so new shared strings land in the
# same relative order the authoring workbook used)
$ws.Range("A8").Value = "kerosene"
$ws.Range("B8").Value = -0.15
$ws.Range("C8").Value = -0.15
$ws.Range("D8").Value = -0.25

$ws.Range("A9").Value = "heavy or residual fuel oil"
$ws.Range("B9").Value = -0.15
$ws.Range("C9").Value = -0.15
$ws.Range("D9").Value = -0.25

$ws.Range("A10").Value = "LPG propane or butane"
$ws.Range("B10").Value = -0.15
$ws.Range("C10").Value = -0.15
$ws.Range("D10").Value = -0.25

$ws.Range("A11").Value = "hydrogen"
$ws.Range("B11").Value = -0.15
$ws.Range("C11").Value = -0.15
$ws.Range("D11").Value = -0.25

# --- Header row (row 1) ------------------------------------------------
$ws.Range("A1").Value = "Elasticity by Fuel (dimensionless)"
$ws.Range("A1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 30

# --- Existing fuel rows (2-6): bump up the row height slightly ---------
$ws.Rows.Item(2).RowHeight = 14.45
$ws.Rows.Item(3).RowHeight = 14.45
$ws.Rows.Item(4).RowHeight = 14.45
$ws.Rows.Item(5).RowHeight = 14.45
$ws.Rows.Item(6).RowHeight = 14.45

# --- Column width tweaks -------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 24.166666666666668
$ws.Range("B1:C1").EntireColumn.ColumnWidth = 19.022135416666668
$ws.Columns.Item(4).ColumnWidth = 13.307291666666666
